$d = $word.ActiveDocument

# Find the paragraph holding the last existing to-do item ("When paying
# utilities, ...") and insert a brand-new to-do item right after it, using
# the same Times New Roman / 12pt (sz 24) formatting as the rest of the list.

$anchorText = "When paying utilities, you don"
$newText = "If someone is visiting jail, and someone goes to jail, it creates an infinite loop of board printing"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($anchorText)) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph containing '$anchorText'"
}

# Insert a new paragraph right after the anchor paragraph.
$target.Range.InsertParagraphAfter()

# The newly created paragraph now sits right after $target.
$newPara = $target.Next()
$newRange = $newPara.Range

# Make sure the paragraph mark itself carries the same run formatting as the
# rest of the list (Times New Roman, sz 24 / szCs 24), then set the text.
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 12

$newRange.Text = $newText
$newRange.Font.Name = "Times New Roman"
$newRange.Font.Size = 12
